$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "bla with text"
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 0.232
$ws.Range("B4").Value = 0.01
$ws.Range("B5").Value = 0.1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = "text"

$ws.Columns("B:B").NumberFormat = $ws.Range("A1").NumberFormat

$ws.Range("B2").Select()
